$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = 'Transferência'
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '29/08/2024'
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = 'TRANSF.REALIZADA PIX SICOOB-MESMA TIT. FAV.: MMG LOCACAO DE MAQUINAS LTDA Transferência Pix MMG LOCACAO DE MAQUINAS LTDA 44.388.803 0001-30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = 3040
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "0.00"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = ""
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = 'MMG LOCACAO DE MAQUINAS LTDA'
$ws.Range("G2").Style = "Normal"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = 'Recebimentos'
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '28/08/2024'
$ws.Range("B3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = 3042.62
$ws.Range("E3").Style = "Normal"
$ws.Range("E3").NumberFormat = "0.00"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = 'Locação de equipamentos'
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '-'
$ws.Range("G3").Style = "Normal"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = 'Despesas Fixas'
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '26/08/2024'
$ws.Range("B4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = 'Contrato Natália Lopez'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 2000
$ws.Range("E4").Style = "Normal"
$ws.Range("E4").NumberFormat = "0.00"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = 'Despesas administrativas'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '-'
$ws.Range("G4").Style = "Normal"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = 'Despesas Fixas'
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '26/08/2024'
$ws.Range("B5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = 'Contrato Natália Lopez'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = 1064.7
$ws.Range("E5").Style = "Normal"
$ws.Range("E5").NumberFormat = "0.00"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = 'Despesas administrativas'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '-'
$ws.Range("G5").Style = "Normal"

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = 'Despesas Variaveis'
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '22/08/2024'
$ws.Range("B6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = 'Locação de equipamentos'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = 1250
$ws.Range("E6").Style = "Normal"
$ws.Range("E6").NumberFormat = "0.00"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = 'Custos operacionais'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = 'MMG LOCACAO DE MAQUINAS LTDA'
$ws.Range("G6").Style = "Normal"

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = 'Recebimentos'
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '22/08/2024'
$ws.Range("B7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = 949.79
$ws.Range("E7").Style = "Normal"
$ws.Range("E7").NumberFormat = "0.00"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = 'Locação de equipamentos'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '-'
$ws.Range("G7").Style = "Normal"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = 'Recebimentos'
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '22/08/2024'
$ws.Range("B8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = 'CR ANTECIPAÇÃO VISA SIPAG_Ant._Visa'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = 447.35
$ws.Range("E8").Style = "Normal"
$ws.Range("E8").NumberFormat = "0.00"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = 'Locação de equipamentos'
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '-'
$ws.Range("G8").Style = "Normal"

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = 'Despesas Variaveis'
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '21/08/2024'
$ws.Range("B9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = 'Alimentação e hospedagem'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = 333
$ws.Range("E9").Style = "Normal"
$ws.Range("E9").NumberFormat = "0.00"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = 'Custos operacionais'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = 'ZOOP TECNOLOGIA & INSTITUICAO DE PAGAMENTO S.A.'
$ws.Range("G9").Style = "Normal"

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = 'Despesas Variaveis'
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '21/08/2024'
$ws.Range("B10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = 'Combustível'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = 1000
$ws.Range("E10").Style = "Normal"
$ws.Range("E10").NumberFormat = "0.00"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = 'Custos operacionais'
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = 'TICKET SOLUCOES HDFGT S/A'
$ws.Range("G10").Style = "Normal"

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = 'Despesas Fixas'
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '21/08/2024'
$ws.Range("B11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = 'Bruna'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = 10500
$ws.Range("E11").Style = "Normal"
$ws.Range("E11").NumberFormat = "0.00"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = 'Empréstimo'
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = 'Bruna Zavati Zavitoski'
$ws.Range("G11").Style = "Normal"

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = 'Recebimentos'
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '21/08/2024'
$ws.Range("B12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = 10750.25
$ws.Range("E12").Style = "Normal"
$ws.Range("E12").NumberFormat = "0.00"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = 'Locação de equipamentos'
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '-'
$ws.Range("G12").Style = "Normal"

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = 'Despesas Fixas'
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '20/08/2024'
$ws.Range("B13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = 'Seguro de equipamento'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = 704.52
$ws.Range("E13").Style = "Normal"
$ws.Range("E13").NumberFormat = "0.00"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = 'Custos operacionais'
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = 'EASY SEGUROS MEDICOS LTDA'
$ws.Range("G13").Style = "Normal"

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = 'Despesas Variaveis'
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '16/08/2024'
$ws.Range("B14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = 'Alimentação e hospedagem'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = 200
$ws.Range("E14").Style = "Normal"
$ws.Range("E14").NumberFormat = "0.00"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = 'Custos operacionais'
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = 'ZOOP TECNOLOGIA & INSTITUICAO DE PAGAMENTO S.A.'
$ws.Range("G14").Style = "Normal"

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = 'Despesas Fixas'
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '16/08/2024'
$ws.Range("B15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = 'Seguro Veicular'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = 227.42
$ws.Range("E15").Style = "Normal"
$ws.Range("E15").NumberFormat = "0.00"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = 'Custos operacionais'
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '-'
$ws.Range("G15").Style = "Normal"

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = 'Despesas Variaveis'
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '16/08/2024'
$ws.Range("B16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = 'Manutenção de veículo'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = 1070
$ws.Range("E16").Style = "Normal"
$ws.Range("E16").NumberFormat = "0.00"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = 'Manutenção de ativos'
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = 'MICHEL ALEXANDRE DE PAULA 39977974802'
$ws.Range("G16").Style = "Normal"

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = 'Despesas Fixas'
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '16/08/2024'
$ws.Range("B17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = 'Bruna'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = 12239.5
$ws.Range("E17").Style = "Normal"
$ws.Range("E17").NumberFormat = "0.00"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = 'Empréstimo'
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = 'Bruna Zavati Zavitoski'
$ws.Range("G17").Style = "Normal"

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = 'Pessoas'
$ws.Range("A18").Style = "Normal"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '16/08/2024'
$ws.Range("B18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = 'Pró-Labore - Guilherme'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = 5131.42
$ws.Range("E18").Style = "Normal"
$ws.Range("E18").NumberFormat = "0.00"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = 'Despesas administrativas'
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = 'Guilherme Borelli'
$ws.Range("G18").Style = "Normal"

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = 'Recebimentos'
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '16/08/2024'
$ws.Range("B19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = 'CR ANTECIPAÇÃO VISA SIPAG_Ant._Visa'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = 9854.44
$ws.Range("E19").Style = "Normal"
$ws.Range("E19").NumberFormat = "0.00"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = 'Locação de equipamentos'
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '-'
$ws.Range("G19").Style = "Normal"

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = 'Recebimentos'
$ws.Range("A20").Style = "Normal"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '16/08/2024'
$ws.Range("B20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = 9854.44
$ws.Range("E20").Style = "Normal"
$ws.Range("E20").NumberFormat = "0.00"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = 'Locação de equipamentos'
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '-'
$ws.Range("G20").Style = "Normal"

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = 'Despesas Variaveis'
$ws.Range("A21").Style = "Normal"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '15/08/2024'
$ws.Range("B21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = 'Locação de equipamentos'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = 2440
$ws.Range("E21").Style = "Normal"
$ws.Range("E21").NumberFormat = "0.00"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = 'Custos operacionais'
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = 'LUME LASER - EQUIPAMENTOS PARA ESTETICA LTDA'
$ws.Range("G21").Style = "Normal"

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = 'Despesas Fixas'
$ws.Range("A22").Style = "Normal"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '15/08/2024'
$ws.Range("B22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = 'Eugênio - Serviços evento'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = 1900
$ws.Range("E22").Style = "Normal"
$ws.Range("E22").NumberFormat = "0.00"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = 'Despesas administrativas'
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '-'
$ws.Range("G22").Style = "Normal"

# Row 23
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = 'Despesas Fixas'
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '14/08/2024'
$ws.Range("B23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = 'Renegociação boletos em atraso'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = 24377.8
$ws.Range("E23").Style = "Normal"
$ws.Range("E23").NumberFormat = "0.00"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = 'Compra de ativos'
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = 'PERSONALITE RECUPERACAO DE CREDITO LTDA'
$ws.Range("G23").Style = "Normal"

# Row 24
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = 'Despesas Variaveis'
$ws.Range("A24").Style = "Normal"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '14/08/2024'
$ws.Range("B24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = 'Locação de equipamentos'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = 3200
$ws.Range("E24").Style = "Normal"
$ws.Range("E24").NumberFormat = "0.00"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = 'Custos operacionais'
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '-'
$ws.Range("G24").Style = "Normal"

# Row 25
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = 'Pessoas'
$ws.Range("A25").Style = "Normal"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '13/08/2024'
$ws.Range("B25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = 'Bruna'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = 6880
$ws.Range("E25").Style = "Normal"
$ws.Range("E25").NumberFormat = "0.00"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = 'Despesas administrativas'
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = 'Bruna Zavati Zavitoski'
$ws.Range("G25").Style = "Normal"

# Row 26
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = 'Recebimentos'
$ws.Range("A26").Style = "Normal"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '13/08/2024'
$ws.Range("B26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = 'CR ANTECIPAÇÃO VISA SIPAG_Ant._Visa'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = 19352.99
$ws.Range("E26").Style = "Normal"
$ws.Range("E26").NumberFormat = "0.00"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = 'Locação de equipamentos'
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '-'
$ws.Range("G26").Style = "Normal"

# Row 27
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = 'Recebimentos'
$ws.Range("A27").Style = "Normal"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '13/08/2024'
$ws.Range("B27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = 23293.38
$ws.Range("E27").Style = "Normal"
$ws.Range("E27").NumberFormat = "0.00"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = 'Locação de equipamentos'
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '-'
$ws.Range("G27").Style = "Normal"

# Row 28
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = 'Despesas Fixas'
$ws.Range("A28").Style = "Normal"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = '12/08/2024'
$ws.Range("B28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = 'Compras evento'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = 1900
$ws.Range("E28").Style = "Normal"
$ws.Range("E28").NumberFormat = "0.00"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = 'Despesas administrativas'
$ws.Range("F28").Style = "Normal"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '-'
$ws.Range("G28").Style = "Normal"

# Row 29
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = 'Recebimentos'
$ws.Range("A29").Style = "Normal"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = '09/08/2024'
$ws.Range("B29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = 'CRED.TRANSF.CONTAS INTERCREDIS REM.: SUELLEN D. DA S. MENDES ESTETICA Transferência Pix SUELLEN D. DA S. MENDES ESTETICA 54.668.555 0001-35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = 497
$ws.Range("E29").Style = "Normal"
$ws.Range("E29").NumberFormat = "0.00"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = 'Locação de equipamentos'
$ws.Range("F29").Style = "Normal"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = 'SUELLEN D. DA S. MENDES ESTETICA'
$ws.Range("G29").Style = "Normal"

# Row 30
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = 'Pessoas'
$ws.Range("A30").Style = "Normal"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = '08/08/2024'
$ws.Range("B30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = 'Pró-Labore - Thauana'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = 1500
$ws.Range("E30").Style = "Normal"
$ws.Range("E30").NumberFormat = "0.00"
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = 'Despesas administrativas'
$ws.Range("F30").Style = "Normal"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '-'
$ws.Range("G30").Style = "Normal"

# Row 31
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = 'Despesas Fixas'
$ws.Range("A31").Style = "Normal"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = '08/08/2024'
$ws.Range("B31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = 'Bruna'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = 8063.51
$ws.Range("E31").Style = "Normal"
$ws.Range("E31").NumberFormat = "0.00"
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = 'Empréstimo'
$ws.Range("F31").Style = "Normal"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = 'Bruna Zavati Zavitoski'
$ws.Range("G31").Style = "Normal"

# Row 32
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = 'Recebimentos'
$ws.Range("A32").Style = "Normal"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = '08/08/2024'
$ws.Range("B32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = 8063.15
$ws.Range("E32").Style = "Normal"
$ws.Range("E32").NumberFormat = "0.00"
$ws.Range("F32").NumberFormat = "@"
$ws.Range("F32").Value = 'Locação de equipamentos'
$ws.Range("F32").Style = "Normal"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '-'
$ws.Range("G32").Style = "Normal"

# Row 33
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = 'Recebimentos'
$ws.Range("A33").Style = "Normal"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = '06/08/2024'
$ws.Range("B33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = 'Ingresso evento'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = 497
$ws.Range("E33").Style = "Normal"
$ws.Range("E33").NumberFormat = "0.00"
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = 'Marketing'
$ws.Range("F33").Style = "Normal"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '-'
$ws.Range("G33").Style = "Normal"

# Row 34
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = 'Pessoas'
$ws.Range("A34").Style = "Normal"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = '05/08/2024'
$ws.Range("B34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = 'Pró-Labore - Guilherme'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = 3000
$ws.Range("E34").Style = "Normal"
$ws.Range("E34").NumberFormat = "0.00"
$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = 'Despesas administrativas'
$ws.Range("F34").Style = "Normal"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = 'Guilherme Borelli'
$ws.Range("G34").Style = "Normal"

# Row 35
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = 'Recebimentos'
$ws.Range("A35").Style = "Normal"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = '05/08/2024'
$ws.Range("B35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = 472.5
$ws.Range("E35").Style = "Normal"
$ws.Range("E35").NumberFormat = "0.00"
$ws.Range("F35").NumberFormat = "@"
$ws.Range("F35").Value = 'Locação de equipamentos'
$ws.Range("F35").Style = "Normal"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '-'
$ws.Range("G35").Style = "Normal"

# Row 36
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = 'Recebimentos'
$ws.Range("A36").Style = "Normal"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = '02/08/2024'
$ws.Range("B36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = 'PIX RECEBIDO - OUTRA IF Recebimento Pix CAROLINE ALMEIDA LEVORCI ***.523.890-**'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = 1000
$ws.Range("E36").Style = "Normal"
$ws.Range("E36").NumberFormat = "0.00"
$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = 'Locação de equipamentos'
$ws.Range("F36").Style = "Normal"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '-'
$ws.Range("G36").Style = "Normal"

# Row 37
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = 'Despesas Fixas'
$ws.Range("A37").Style = "Normal"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = '01/08/2024'
$ws.Range("B37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = 'Bradesco fiananciamento - Pixie'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = 7809.48
$ws.Range("E37").Style = "Normal"
$ws.Range("E37").NumberFormat = "0.00"
$ws.Range("F37").NumberFormat = "@"
$ws.Range("F37").Value = 'Compra de ativos'
$ws.Range("F37").Style = "Normal"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '-'
$ws.Range("G37").Style = "Normal"

# Row 38
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = 'Recebimentos'
$ws.Range("A38").Style = "Normal"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = '01/08/2024'
$ws.Range("B38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = 'Ingresso evento'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = 497
$ws.Range("E38").Style = "Normal"
$ws.Range("E38").NumberFormat = "0.00"
$ws.Range("F38").NumberFormat = "@"
$ws.Range("F38").Value = 'Marketing'
$ws.Range("F38").Style = "Normal"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '-'
$ws.Range("G38").Style = "Normal"

# Row 39
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = 'Recebimentos'
$ws.Range("A39").Style = "Normal"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = '01/08/2024'
$ws.Range("B39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = 'Ingresso evento'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = 497
$ws.Range("E39").Style = "Normal"
$ws.Range("E39").NumberFormat = "0.00"
$ws.Range("F39").NumberFormat = "@"
$ws.Range("F39").Value = 'Marketing'
$ws.Range("F39").Style = "Normal"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = 'INFINITY COMPANY BY GG LTDA'
$ws.Range("G39").Style = "Normal"

# Row 40
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = 'Recebimentos'
$ws.Range("A40").Style = "Normal"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = '01/08/2024'
$ws.Range("B40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = 'CR ANTECIPAÇÃO MASTERCARD SIPAG_Ant._Mastercard'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = 2729.55
$ws.Range("E40").Style = "Normal"
$ws.Range("E40").NumberFormat = "0.00"
$ws.Range("F40").NumberFormat = "@"
$ws.Range("F40").Value = 'Locação de equipamentos'
$ws.Range("F40").Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '-'
$ws.Range("G40").Style = "Normal"

# Row 24 extra columns
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = ""
$ws.Range("C24").Style = "Normal"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = 'Sim'
$ws.Range("H24").Style = "Normal"
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = 'Preenchido Automaticamente'
$ws.Range("I24").Style = "Normal"
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = 'Sicoob 2'
$ws.Range("J24").Style = "Normal"
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = ""
$ws.Range("K24").Style = "Normal"
$ws.Range("L24").NumberFormat = "@"
$ws.Range("L24").Value = ""
$ws.Range("L24").Style = "Normal"
$ws.Range("M24").NumberFormat = "@"
$ws.Range("M24").Value = 'Indefinido'
$ws.Range("M24").Style = "Normal"
$ws.Range("N24").NumberFormat = "@"
$ws.Range("N24").Value = ""
$ws.Range("N24").Style = "Normal"
$ws.Range("O24").NumberFormat = "@"
$ws.Range("O24").Value = ""
$ws.Range("O24").Style = "Normal"
$ws.Range("P24").NumberFormat = "@"
$ws.Range("P24").Value = ""
$ws.Range("P24").Style = "Normal"
$ws.Range("Q24").NumberFormat = "@"
$ws.Range("Q24").Value = ""
$ws.Range("Q24").Style = "Normal"
$ws.Range("R24").NumberFormat = "@"
$ws.Range("R24").Value = ""
$ws.Range("R24").Style = "Normal"
$ws.Range("S24").NumberFormat = "@"
$ws.Range("S24").Value = ""
$ws.Range("S24").Style = "Normal"

# Row 25 extra columns
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = ""
$ws.Range("C25").Style = "Normal"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = 'Sim'
$ws.Range("H25").Style = "Normal"
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I25").Value = 'Preenchido Automaticamente'
$ws.Range("I25").Style = "Normal"
$ws.Range("J25").NumberFormat = "@"
$ws.Range("J25").Value = 'Sicoob 2'
$ws.Range("J25").Style = "Normal"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = ""
$ws.Range("K25").Style = "Normal"
$ws.Range("L25").NumberFormat = "@"
$ws.Range("L25").Value = ""
$ws.Range("L25").Style = "Normal"
$ws.Range("M25").NumberFormat = "@"
$ws.Range("M25").Value = 'Indefinido'
$ws.Range("M25").Style = "Normal"
$ws.Range("N25").NumberFormat = "@"
$ws.Range("N25").Value = ""
$ws.Range("N25").Style = "Normal"
$ws.Range("O25").NumberFormat = "@"
$ws.Range("O25").Value = ""
$ws.Range("O25").Style = "Normal"
$ws.Range("P25").NumberFormat = "@"
$ws.Range("P25").Value = ""
$ws.Range("P25").Style = "Normal"
$ws.Range("Q25").NumberFormat = "@"
$ws.Range("Q25").Value = ""
$ws.Range("Q25").Style = "Normal"
$ws.Range("R25").NumberFormat = "@"
$ws.Range("R25").Value = ""
$ws.Range("R25").Style = "Normal"
$ws.Range("S25").NumberFormat = "@"
$ws.Range("S25").Value = ""
$ws.Range("S25").Style = "Normal"

# Row 26 extra columns
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = ""
$ws.Range("C26").Style = "Normal"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = 'Sim'
$ws.Range("H26").Style = "Normal"
$ws.Range("I26").NumberFormat = "@"
$ws.Range("I26").Value = 'Preenchido Automaticamente'
$ws.Range("I26").Style = "Normal"
$ws.Range("J26").NumberFormat = "@"
$ws.Range("J26").Value = 'Sicoob 2'
$ws.Range("J26").Style = "Normal"
$ws.Range("K26").NumberFormat = "@"
$ws.Range("K26").Value = ""
$ws.Range("K26").Style = "Normal"
$ws.Range("L26").NumberFormat = "@"
$ws.Range("L26").Value = ""
$ws.Range("L26").Style = "Normal"
$ws.Range("M26").NumberFormat = "@"
$ws.Range("M26").Value = 'Indefinido'
$ws.Range("M26").Style = "Normal"
$ws.Range("N26").NumberFormat = "@"
$ws.Range("N26").Value = ""
$ws.Range("N26").Style = "Normal"
$ws.Range("O26").NumberFormat = "@"
$ws.Range("O26").Value = ""
$ws.Range("O26").Style = "Normal"
$ws.Range("P26").NumberFormat = "@"
$ws.Range("P26").Value = ""
$ws.Range("P26").Style = "Normal"
$ws.Range("Q26").NumberFormat = "@"
$ws.Range("Q26").Value = ""
$ws.Range("Q26").Style = "Normal"
$ws.Range("R26").NumberFormat = "@"
$ws.Range("R26").Value = ""
$ws.Range("R26").Style = "Normal"
$ws.Range("S26").NumberFormat = "@"
$ws.Range("S26").Value = ""
$ws.Range("S26").Style = "Normal"

# Row 27 extra columns
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = ""
$ws.Range("C27").Style = "Normal"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = 'Sim'
$ws.Range("H27").Style = "Normal"
$ws.Range("I27").NumberFormat = "@"
$ws.Range("I27").Value = 'Preenchido Automaticamente'
$ws.Range("I27").Style = "Normal"
$ws.Range("J27").NumberFormat = "@"
$ws.Range("J27").Value = 'Sicoob 2'
$ws.Range("J27").Style = "Normal"
$ws.Range("K27").NumberFormat = "@"
$ws.Range("K27").Value = ""
$ws.Range("K27").Style = "Normal"
$ws.Range("L27").NumberFormat = "@"
$ws.Range("L27").Value = ""
$ws.Range("L27").Style = "Normal"
$ws.Range("M27").NumberFormat = "@"
$ws.Range("M27").Value = 'Indefinido'
$ws.Range("M27").Style = "Normal"
$ws.Range("N27").NumberFormat = "@"
$ws.Range("N27").Value = ""
$ws.Range("N27").Style = "Normal"
$ws.Range("O27").NumberFormat = "@"
$ws.Range("O27").Value = ""
$ws.Range("O27").Style = "Normal"
$ws.Range("P27").NumberFormat = "@"
$ws.Range("P27").Value = ""
$ws.Range("P27").Style = "Normal"
$ws.Range("Q27").NumberFormat = "@"
$ws.Range("Q27").Value = ""
$ws.Range("Q27").Style = "Normal"
$ws.Range("R27").NumberFormat = "@"
$ws.Range("R27").Value = ""
$ws.Range("R27").Style = "Normal"
$ws.Range("S27").NumberFormat = "@"
$ws.Range("S27").Value = ""
$ws.Range("S27").Style = "Normal"

# Row 28 extra columns
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = ""
$ws.Range("C28").Style = "Normal"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = 'Sim'
$ws.Range("H28").Style = "Normal"
$ws.Range("I28").NumberFormat = "@"
$ws.Range("I28").Value = 'Preenchido Automaticamente'
$ws.Range("I28").Style = "Normal"
$ws.Range("J28").NumberFormat = "@"
$ws.Range("J28").Value = 'Sicoob 2'
$ws.Range("J28").Style = "Normal"
$ws.Range("K28").NumberFormat = "@"
$ws.Range("K28").Value = ""
$ws.Range("K28").Style = "Normal"
$ws.Range("L28").NumberFormat = "@"
$ws.Range("L28").Value = ""
$ws.Range("L28").Style = "Normal"
$ws.Range("M28").NumberFormat = "@"
$ws.Range("M28").Value = 'Indefinido'
$ws.Range("M28").Style = "Normal"
$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value = ""
$ws.Range("N28").Style = "Normal"
$ws.Range("O28").NumberFormat = "@"
$ws.Range("O28").Value = ""
$ws.Range("O28").Style = "Normal"
$ws.Range("P28").NumberFormat = "@"
$ws.Range("P28").Value = ""
$ws.Range("P28").Style = "Normal"
$ws.Range("Q28").NumberFormat = "@"
$ws.Range("Q28").Value = ""
$ws.Range("Q28").Style = "Normal"
$ws.Range("R28").NumberFormat = "@"
$ws.Range("R28").Value = ""
$ws.Range("R28").Style = "Normal"
$ws.Range("S28").NumberFormat = "@"
$ws.Range("S28").Value = ""
$ws.Range("S28").Style = "Normal"

# Row 29 extra columns
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = ""
$ws.Range("C29").Style = "Normal"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = 'Sim'
$ws.Range("H29").Style = "Normal"
$ws.Range("I29").NumberFormat = "@"
$ws.Range("I29").Value = 'Preenchido Automaticamente'
$ws.Range("I29").Style = "Normal"
$ws.Range("J29").NumberFormat = "@"
$ws.Range("J29").Value = 'Sicoob 2'
$ws.Range("J29").Style = "Normal"
$ws.Range("K29").NumberFormat = "@"
$ws.Range("K29").Value = ""
$ws.Range("K29").Style = "Normal"
$ws.Range("L29").NumberFormat = "@"
$ws.Range("L29").Value = ""
$ws.Range("L29").Style = "Normal"
$ws.Range("M29").NumberFormat = "@"
$ws.Range("M29").Value = 'Indefinido'
$ws.Range("M29").Style = "Normal"
$ws.Range("N29").NumberFormat = "@"
$ws.Range("N29").Value = ""
$ws.Range("N29").Style = "Normal"
$ws.Range("O29").NumberFormat = "@"
$ws.Range("O29").Value = ""
$ws.Range("O29").Style = "Normal"
$ws.Range("P29").NumberFormat = "@"
$ws.Range("P29").Value = ""
$ws.Range("P29").Style = "Normal"
$ws.Range("Q29").NumberFormat = "@"
$ws.Range("Q29").Value = ""
$ws.Range("Q29").Style = "Normal"
$ws.Range("R29").NumberFormat = "@"
$ws.Range("R29").Value = ""
$ws.Range("R29").Style = "Normal"
$ws.Range("S29").NumberFormat = "@"
$ws.Range("S29").Value = ""
$ws.Range("S29").Style = "Normal"

# Row 30 extra columns
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = ""
$ws.Range("C30").Style = "Normal"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = 'Sim'
$ws.Range("H30").Style = "Normal"
$ws.Range("I30").NumberFormat = "@"
$ws.Range("I30").Value = 'Preenchido Automaticamente'
$ws.Range("I30").Style = "Normal"
$ws.Range("J30").NumberFormat = "@"
$ws.Range("J30").Value = 'Sicoob 2'
$ws.Range("J30").Style = "Normal"
$ws.Range("K30").NumberFormat = "@"
$ws.Range("K30").Value = ""
$ws.Range("K30").Style = "Normal"
$ws.Range("L30").NumberFormat = "@"
$ws.Range("L30").Value = ""
$ws.Range("L30").Style = "Normal"
$ws.Range("M30").NumberFormat = "@"
$ws.Range("M30").Value = 'Indefinido'
$ws.Range("M30").Style = "Normal"
$ws.Range("N30").NumberFormat = "@"
$ws.Range("N30").Value = ""
$ws.Range("N30").Style = "Normal"
$ws.Range("O30").NumberFormat = "@"
$ws.Range("O30").Value = ""
$ws.Range("O30").Style = "Normal"
$ws.Range("P30").NumberFormat = "@"
$ws.Range("P30").Value = ""
$ws.Range("P30").Style = "Normal"
$ws.Range("Q30").NumberFormat = "@"
$ws.Range("Q30").Value = ""
$ws.Range("Q30").Style = "Normal"
$ws.Range("R30").NumberFormat = "@"
$ws.Range("R30").Value = ""
$ws.Range("R30").Style = "Normal"
$ws.Range("S30").NumberFormat = "@"
$ws.Range("S30").Value = ""
$ws.Range("S30").Style = "Normal"

# Row 31 extra columns
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = ""
$ws.Range("C31").Style = "Normal"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = 'Sim'
$ws.Range("H31").Style = "Normal"
$ws.Range("I31").NumberFormat = "@"
$ws.Range("I31").Value = 'Preenchido Automaticamente'
$ws.Range("I31").Style = "Normal"
$ws.Range("J31").NumberFormat = "@"
$ws.Range("J31").Value = 'Sicoob 2'
$ws.Range("J31").Style = "Normal"
$ws.Range("K31").NumberFormat = "@"
$ws.Range("K31").Value = ""
$ws.Range("K31").Style = "Normal"
$ws.Range("L31").NumberFormat = "@"
$ws.Range("L31").Value = ""
$ws.Range("L31").Style = "Normal"
$ws.Range("M31").NumberFormat = "@"
$ws.Range("M31").Value = 'Indefinido'
$ws.Range("M31").Style = "Normal"
$ws.Range("N31").NumberFormat = "@"
$ws.Range("N31").Value = ""
$ws.Range("N31").Style = "Normal"
$ws.Range("O31").NumberFormat = "@"
$ws.Range("O31").Value = ""
$ws.Range("O31").Style = "Normal"
$ws.Range("P31").NumberFormat = "@"
$ws.Range("P31").Value = ""
$ws.Range("P31").Style = "Normal"
$ws.Range("Q31").NumberFormat = "@"
$ws.Range("Q31").Value = ""
$ws.Range("Q31").Style = "Normal"
$ws.Range("R31").NumberFormat = "@"
$ws.Range("R31").Value = ""
$ws.Range("R31").Style = "Normal"
$ws.Range("S31").NumberFormat = "@"
$ws.Range("S31").Value = ""
$ws.Range("S31").Style = "Normal"

# Row 32 extra columns
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = ""
$ws.Range("C32").Style = "Normal"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = 'Sim'
$ws.Range("H32").Style = "Normal"
$ws.Range("I32").NumberFormat = "@"
$ws.Range("I32").Value = 'Preenchido Automaticamente'
$ws.Range("I32").Style = "Normal"
$ws.Range("J32").NumberFormat = "@"
$ws.Range("J32").Value = 'Sicoob 2'
$ws.Range("J32").Style = "Normal"
$ws.Range("K32").NumberFormat = "@"
$ws.Range("K32").Value = ""
$ws.Range("K32").Style = "Normal"
$ws.Range("L32").NumberFormat = "@"
$ws.Range("L32").Value = ""
$ws.Range("L32").Style = "Normal"
$ws.Range("M32").NumberFormat = "@"
$ws.Range("M32").Value = 'Indefinido'
$ws.Range("M32").Style = "Normal"
$ws.Range("N32").NumberFormat = "@"
$ws.Range("N32").Value = ""
$ws.Range("N32").Style = "Normal"
$ws.Range("O32").NumberFormat = "@"
$ws.Range("O32").Value = ""
$ws.Range("O32").Style = "Normal"
$ws.Range("P32").NumberFormat = "@"
$ws.Range("P32").Value = ""
$ws.Range("P32").Style = "Normal"
$ws.Range("Q32").NumberFormat = "@"
$ws.Range("Q32").Value = ""
$ws.Range("Q32").Style = "Normal"
$ws.Range("R32").NumberFormat = "@"
$ws.Range("R32").Value = ""
$ws.Range("R32").Style = "Normal"
$ws.Range("S32").NumberFormat = "@"
$ws.Range("S32").Value = ""
$ws.Range("S32").Style = "Normal"

# Row 33 extra columns
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = ""
$ws.Range("C33").Style = "Normal"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = 'Sim'
$ws.Range("H33").Style = "Normal"
$ws.Range("I33").NumberFormat = "@"
$ws.Range("I33").Value = 'Preenchido Automaticamente'
$ws.Range("I33").Style = "Normal"
$ws.Range("J33").NumberFormat = "@"
$ws.Range("J33").Value = 'Sicoob 2'
$ws.Range("J33").Style = "Normal"
$ws.Range("K33").NumberFormat = "@"
$ws.Range("K33").Value = ""
$ws.Range("K33").Style = "Normal"
$ws.Range("L33").NumberFormat = "@"
$ws.Range("L33").Value = ""
$ws.Range("L33").Style = "Normal"
$ws.Range("M33").NumberFormat = "@"
$ws.Range("M33").Value = 'Indefinido'
$ws.Range("M33").Style = "Normal"
$ws.Range("N33").NumberFormat = "@"
$ws.Range("N33").Value = ""
$ws.Range("N33").Style = "Normal"
$ws.Range("O33").NumberFormat = "@"
$ws.Range("O33").Value = ""
$ws.Range("O33").Style = "Normal"
$ws.Range("P33").NumberFormat = "@"
$ws.Range("P33").Value = ""
$ws.Range("P33").Style = "Normal"
$ws.Range("Q33").NumberFormat = "@"
$ws.Range("Q33").Value = ""
$ws.Range("Q33").Style = "Normal"
$ws.Range("R33").NumberFormat = "@"
$ws.Range("R33").Value = ""
$ws.Range("R33").Style = "Normal"
$ws.Range("S33").NumberFormat = "@"
$ws.Range("S33").Value = ""
$ws.Range("S33").Style = "Normal"

# Row 34 extra columns
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = ""
$ws.Range("C34").Style = "Normal"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = 'Sim'
$ws.Range("H34").Style = "Normal"
$ws.Range("I34").NumberFormat = "@"
$ws.Range("I34").Value = 'Preenchido Automaticamente'
$ws.Range("I34").Style = "Normal"
$ws.Range("J34").NumberFormat = "@"
$ws.Range("J34").Value = 'Sicoob 2'
$ws.Range("J34").Style = "Normal"
$ws.Range("K34").NumberFormat = "@"
$ws.Range("K34").Value = ""
$ws.Range("K34").Style = "Normal"
$ws.Range("L34").NumberFormat = "@"
$ws.Range("L34").Value = ""
$ws.Range("L34").Style = "Normal"
$ws.Range("M34").NumberFormat = "@"
$ws.Range("M34").Value = 'Indefinido'
$ws.Range("M34").Style = "Normal"
$ws.Range("N34").NumberFormat = "@"
$ws.Range("N34").Value = ""
$ws.Range("N34").Style = "Normal"
$ws.Range("O34").NumberFormat = "@"
$ws.Range("O34").Value = ""
$ws.Range("O34").Style = "Normal"
$ws.Range("P34").NumberFormat = "@"
$ws.Range("P34").Value = ""
$ws.Range("P34").Style = "Normal"
$ws.Range("Q34").NumberFormat = "@"
$ws.Range("Q34").Value = ""
$ws.Range("Q34").Style = "Normal"
$ws.Range("R34").NumberFormat = "@"
$ws.Range("R34").Value = ""
$ws.Range("R34").Style = "Normal"
$ws.Range("S34").NumberFormat = "@"
$ws.Range("S34").Value = ""
$ws.Range("S34").Style = "Normal"

# Row 35 extra columns
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = ""
$ws.Range("C35").Style = "Normal"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = 'Sim'
$ws.Range("H35").Style = "Normal"
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = 'Preenchido Automaticamente'
$ws.Range("I35").Style = "Normal"
$ws.Range("J35").NumberFormat = "@"
$ws.Range("J35").Value = 'Sicoob 2'
$ws.Range("J35").Style = "Normal"
$ws.Range("K35").NumberFormat = "@"
$ws.Range("K35").Value = ""
$ws.Range("K35").Style = "Normal"
$ws.Range("L35").NumberFormat = "@"
$ws.Range("L35").Value = ""
$ws.Range("L35").Style = "Normal"
$ws.Range("M35").NumberFormat = "@"
$ws.Range("M35").Value = 'Indefinido'
$ws.Range("M35").Style = "Normal"
$ws.Range("N35").NumberFormat = "@"
$ws.Range("N35").Value = ""
$ws.Range("N35").Style = "Normal"
$ws.Range("O35").NumberFormat = "@"
$ws.Range("O35").Value = ""
$ws.Range("O35").Style = "Normal"
$ws.Range("P35").NumberFormat = "@"
$ws.Range("P35").Value = ""
$ws.Range("P35").Style = "Normal"
$ws.Range("Q35").NumberFormat = "@"
$ws.Range("Q35").Value = ""
$ws.Range("Q35").Style = "Normal"
$ws.Range("R35").NumberFormat = "@"
$ws.Range("R35").Value = ""
$ws.Range("R35").Style = "Normal"
$ws.Range("S35").NumberFormat = "@"
$ws.Range("S35").Value = ""
$ws.Range("S35").Style = "Normal"

# Row 36 extra columns
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = ""
$ws.Range("C36").Style = "Normal"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = 'Sim'
$ws.Range("H36").Style = "Normal"
$ws.Range("I36").NumberFormat = "@"
$ws.Range("I36").Value = 'Preenchido Automaticamente'
$ws.Range("I36").Style = "Normal"
$ws.Range("J36").NumberFormat = "@"
$ws.Range("J36").Value = 'Sicoob 2'
$ws.Range("J36").Style = "Normal"
$ws.Range("K36").NumberFormat = "@"
$ws.Range("K36").Value = ""
$ws.Range("K36").Style = "Normal"
$ws.Range("L36").NumberFormat = "@"
$ws.Range("L36").Value = ""
$ws.Range("L36").Style = "Normal"
$ws.Range("M36").NumberFormat = "@"
$ws.Range("M36").Value = 'Indefinido'
$ws.Range("M36").Style = "Normal"
$ws.Range("N36").NumberFormat = "@"
$ws.Range("N36").Value = ""
$ws.Range("N36").Style = "Normal"
$ws.Range("O36").NumberFormat = "@"
$ws.Range("O36").Value = ""
$ws.Range("O36").Style = "Normal"
$ws.Range("P36").NumberFormat = "@"
$ws.Range("P36").Value = ""
$ws.Range("P36").Style = "Normal"
$ws.Range("Q36").NumberFormat = "@"
$ws.Range("Q36").Value = ""
$ws.Range("Q36").Style = "Normal"
$ws.Range("R36").NumberFormat = "@"
$ws.Range("R36").Value = ""
$ws.Range("R36").Style = "Normal"
$ws.Range("S36").NumberFormat = "@"
$ws.Range("S36").Value = ""
$ws.Range("S36").Style = "Normal"

# Row 37 extra columns
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = ""
$ws.Range("C37").Style = "Normal"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = 'Sim'
$ws.Range("H37").Style = "Normal"
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = 'Preenchido Automaticamente'
$ws.Range("I37").Style = "Normal"
$ws.Range("J37").NumberFormat = "@"
$ws.Range("J37").Value = 'Sicoob 2'
$ws.Range("J37").Style = "Normal"
$ws.Range("K37").NumberFormat = "@"
$ws.Range("K37").Value = ""
$ws.Range("K37").Style = "Normal"
$ws.Range("L37").NumberFormat = "@"
$ws.Range("L37").Value = ""
$ws.Range("L37").Style = "Normal"
$ws.Range("M37").NumberFormat = "@"
$ws.Range("M37").Value = 'Indefinido'
$ws.Range("M37").Style = "Normal"
$ws.Range("N37").NumberFormat = "@"
$ws.Range("N37").Value = ""
$ws.Range("N37").Style = "Normal"
$ws.Range("O37").NumberFormat = "@"
$ws.Range("O37").Value = ""
$ws.Range("O37").Style = "Normal"
$ws.Range("P37").NumberFormat = "@"
$ws.Range("P37").Value = ""
$ws.Range("P37").Style = "Normal"
$ws.Range("Q37").NumberFormat = "@"
$ws.Range("Q37").Value = ""
$ws.Range("Q37").Style = "Normal"
$ws.Range("R37").NumberFormat = "@"
$ws.Range("R37").Value = ""
$ws.Range("R37").Style = "Normal"
$ws.Range("S37").NumberFormat = "@"
$ws.Range("S37").Value = ""
$ws.Range("S37").Style = "Normal"

# Row 38 extra columns
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = ""
$ws.Range("C38").Style = "Normal"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = 'Sim'
$ws.Range("H38").Style = "Normal"
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = 'Preenchido Automaticamente'
$ws.Range("I38").Style = "Normal"
$ws.Range("J38").NumberFormat = "@"
$ws.Range("J38").Value = 'Sicoob 2'
$ws.Range("J38").Style = "Normal"
$ws.Range("K38").NumberFormat = "@"
$ws.Range("K38").Value = ""
$ws.Range("K38").Style = "Normal"
$ws.Range("L38").NumberFormat = "@"
$ws.Range("L38").Value = ""
$ws.Range("L38").Style = "Normal"
$ws.Range("M38").NumberFormat = "@"
$ws.Range("M38").Value = 'Indefinido'
$ws.Range("M38").Style = "Normal"
$ws.Range("N38").NumberFormat = "@"
$ws.Range("N38").Value = ""
$ws.Range("N38").Style = "Normal"
$ws.Range("O38").NumberFormat = "@"
$ws.Range("O38").Value = ""
$ws.Range("O38").Style = "Normal"
$ws.Range("P38").NumberFormat = "@"
$ws.Range("P38").Value = ""
$ws.Range("P38").Style = "Normal"
$ws.Range("Q38").NumberFormat = "@"
$ws.Range("Q38").Value = ""
$ws.Range("Q38").Style = "Normal"
$ws.Range("R38").NumberFormat = "@"
$ws.Range("R38").Value = ""
$ws.Range("R38").Style = "Normal"
$ws.Range("S38").NumberFormat = "@"
$ws.Range("S38").Value = ""
$ws.Range("S38").Style = "Normal"

# Row 39 extra columns
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = ""
$ws.Range("C39").Style = "Normal"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = 'Sim'
$ws.Range("H39").Style = "Normal"
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = 'Preenchido Automaticamente'
$ws.Range("I39").Style = "Normal"
$ws.Range("J39").NumberFormat = "@"
$ws.Range("J39").Value = 'Sicoob 2'
$ws.Range("J39").Style = "Normal"
$ws.Range("K39").NumberFormat = "@"
$ws.Range("K39").Value = ""
$ws.Range("K39").Style = "Normal"
$ws.Range("L39").NumberFormat = "@"
$ws.Range("L39").Value = ""
$ws.Range("L39").Style = "Normal"
$ws.Range("M39").NumberFormat = "@"
$ws.Range("M39").Value = 'Indefinido'
$ws.Range("M39").Style = "Normal"
$ws.Range("N39").NumberFormat = "@"
$ws.Range("N39").Value = ""
$ws.Range("N39").Style = "Normal"
$ws.Range("O39").NumberFormat = "@"
$ws.Range("O39").Value = ""
$ws.Range("O39").Style = "Normal"
$ws.Range("P39").NumberFormat = "@"
$ws.Range("P39").Value = ""
$ws.Range("P39").Style = "Normal"
$ws.Range("Q39").NumberFormat = "@"
$ws.Range("Q39").Value = ""
$ws.Range("Q39").Style = "Normal"
$ws.Range("R39").NumberFormat = "@"
$ws.Range("R39").Value = ""
$ws.Range("R39").Style = "Normal"
$ws.Range("S39").NumberFormat = "@"
$ws.Range("S39").Value = ""
$ws.Range("S39").Style = "Normal"

# Row 40 extra columns
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = ""
$ws.Range("C40").Style = "Normal"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = 'Sim'
$ws.Range("H40").Style = "Normal"
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value = 'Preenchido Automaticamente'
$ws.Range("I40").Style = "Normal"
$ws.Range("J40").NumberFormat = "@"
$ws.Range("J40").Value = 'Sicoob 2'
$ws.Range("J40").Style = "Normal"
$ws.Range("K40").NumberFormat = "@"
$ws.Range("K40").Value = ""
$ws.Range("K40").Style = "Normal"
$ws.Range("L40").NumberFormat = "@"
$ws.Range("L40").Value = ""
$ws.Range("L40").Style = "Normal"
$ws.Range("M40").NumberFormat = "@"
$ws.Range("M40").Value = 'Indefinido'
$ws.Range("M40").Style = "Normal"
$ws.Range("N40").NumberFormat = "@"
$ws.Range("N40").Value = ""
$ws.Range("N40").Style = "Normal"
$ws.Range("O40").NumberFormat = "@"
$ws.Range("O40").Value = ""
$ws.Range("O40").Style = "Normal"
$ws.Range("P40").NumberFormat = "@"
$ws.Range("P40").Value = ""
$ws.Range("P40").Style = "Normal"
$ws.Range("Q40").NumberFormat = "@"
$ws.Range("Q40").Value = ""
$ws.Range("Q40").Style = "Normal"
$ws.Range("R40").NumberFormat = "@"
$ws.Range("R40").Value = ""
$ws.Range("R40").Style = "Normal"
$ws.Range("S40").NumberFormat = "@"
$ws.Range("S40").Value = ""
$ws.Range("S40").Style = "Normal"
